$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks so we can rebuild them in the new row order
$ws.Hyperlinks.Delete()

# Target data after the edit (title, timestamp, historical distance, time bucket, uri) for rows 2-11
$rows = @(
  @('Hawaii Voter Surveys: How Different Groups Voted', '2020-11-03T19:12:47UTC', '0', 'day_0', 'https://www.nytimes.com/interactive/2020/11/03/us/elections/ap-polls-hawaii.html'),
  @('Presidential Ratings', '1-01-01T00:00:00UTC', 'unknown', 'unknown', 'https://insideelections.com/ratings/president'),
  @('Hawaii Primary Election Results 2020', '2020-05-23T18:15:40UTC', 'unknown', 'unknown', 'https://www.nytimes.com/interactive/2020/05/23/us/elections/results-hawaii-president-democrat-primary-election.html'),
  @('Find Your Local League', '1-01-01T00:00:00UTC', 'unknown', 'unknown', 'https://www.lwv.org/local-leagues/find-local-league'),
  @('Nevada and South Carolina GOP cancel 2020 presidential nominating contests', '2019-09-07T18:09:00UTC', 'unknown', 'unknown', 'https://abcnews.go.com/Politics/trump-gop-canceling-gop-primaries-caucuses/story?id=65436462'),
  @('2020 President - Sabato''s Crystal Ball', '1-01-01T00:00:00UTC', 'unknown', 'unknown', 'http://centerforpolitics.org/crystalball/2020-president/'),
  @('Green Party of Hawai''i Presidential Preference Poll', '1-01-01T00:00:00UTC', 'unknown', 'unknown', 'https://www.opavote.com/results/6550183955398656'),
  @('Biden dominates the electoral map, but here''s how the race could tighten', '2020-08-06T13:13:00UTC', 'unknown', 'unknown', 'https://www.nbcnews.com/politics/meet-the-press/biden-dominates-electoral-map-here-s-how-race-could-tighten-n1236001'),
  @('How Asian Americans Are Thinking About The 2020 Election', '2020-09-18T07:00:02UTC', 'unknown', 'unknown', 'https://fivethirtyeight.com/features/how-asian-americans-are-thinking-about-the-2020-election/'),
  @('2020 Election Forecast', '2020-08-12T06:30:00UTC', 'unknown', 'unknown', 'https://projects.fivethirtyeight.com/2020-election-forecast/')
)

for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = $i + 2
  $row = $rows[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  if ($r -eq 2) {
    $ws.Cells.Item($r, 3).Value = 0
  } else {
    $ws.Cells.Item($r, 3).Value = $row[2]
  }
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Hyperlinks.Add($ws.Cells.Item($r, 5), $row[4])
}

# Restore the Hyperlink style on column E (Hyperlinks.Add can leave an incorrect style)
for ($r = 2; $r -le 11; $r++) {
  $ws.Cells.Item($r, 5).Style = "Hyperlink"
}

Write-Host "Update complete"
